# Re-applies the re-scrape: for several match-days the rows describing the
# individual fixtures (columns F..V - teams, scores, odds, timestamps, URL)
# were re-ordered among themselves (the "Indice"/date columns A..E stayed on
# the same row), and one brand-new fixture (Luton vs Liverpool) was appended
# as the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Destination row -> source row (the content that ends up on the destination
# row is whatever currently sits on the source row). Columns F..V only.
$rowMap = @{
    4=7; 7=4;
    23=24; 24=25; 25=26; 26=23;
    28=29; 29=28;
    44=45; 45=46; 46=44;
    51=53; 53=51;
    56=59; 57=58; 58=57; 59=56;
    62=67; 63=66; 64=62; 65=63; 66=64; 67=65;
    73=76; 74=73; 75=74; 76=75;
    83=84; 84=85; 85=83;
    86=87; 87=86;
    94=95; 95=94;
    99=100; 100=99;
    103=104; 104=105; 105=106; 106=107; 107=103;
}

$firstCol = 6   # column F
$lastCol  = 22  # column V

# Snapshot every row that participates in the shuffle *before* writing
# anything, since several of them form rotation cycles (e.g. 23 <- 24 <- 25
# <- 26 <- 23) where the source for one destination is itself a destination.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values | Sort-Object -Unique) {
    $vals = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $vals[$c] = $ws.Cells.Item($srcRow, $c).Value2
    }
    $snapshot[$srcRow] = $vals
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $vals = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $ws.Cells.Item($destRow, $c).Value = $vals[$c]
    }
}

# Brand-new fixture appended as row 110 (Indice 109): Luton vs Liverpool.
$newRow = 110

# Give the new row's index cell the same style as the other "Indice" cells,
# and the date cell the same style as the other "data_partida" cells, so the
# appended row matches the rest of the table's formatting. Copy() (unlike
# assigning .Style from another cell) reuses the existing style record
# instead of minting a new one. The values get overwritten right after.
$ws.Cells.Item($newRow - 1, 1).Copy($ws.Cells.Item($newRow, 1))
$ws.Cells.Item($newRow - 1, 5).Copy($ws.Cells.Item($newRow, 5))

$ws.Cells.Item($newRow, 1).Value  = 109
$ws.Cells.Item($newRow, 2).Value  = "england"
$ws.Cells.Item($newRow, 3).Value  = "premier-league"
$ws.Cells.Item($newRow, 4).Value  = "2023-2024"
$ws.Cells.Item($newRow, 5).Value  = 45235.72916666666
$ws.Cells.Item($newRow, 6).Value  = "Luton"
$ws.Cells.Item($newRow, 7).Value  = 1
$ws.Cells.Item($newRow, 8).Value  = "Liverpool"
$ws.Cells.Item($newRow, 9).Value  = 1
$ws.Cells.Item($newRow, 10).Value = 7.79
$ws.Cells.Item($newRow, 11).Value = "23/10/2023 15:49"
$ws.Cells.Item($newRow, 12).Value = 15
$ws.Cells.Item($newRow, 13).Value = "05/11/2023 17:28"
$ws.Cells.Item($newRow, 14).Value = 5.75
$ws.Cells.Item($newRow, 15).Value = "23/10/2023 15:49"
$ws.Cells.Item($newRow, 16).Value = 8.5
$ws.Cells.Item($newRow, 17).Value = "05/11/2023 17:28"
$ws.Cells.Item($newRow, 18).Value = 1.31
$ws.Cells.Item($newRow, 19).Value = "23/10/2023 15:49"
$ws.Cells.Item($newRow, 20).Value = 1.18
$ws.Cells.Item($newRow, 21).Value = "05/11/2023 17:27"
$ws.Cells.Item($newRow, 22).Value = "https://www.betexplorer.com/football/england/premier-league/luton-liverpool/W0vgcfWu/"

Write-Output "done"
